$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: Z1, AA1 (copy format from Y1, then set text) ---
$ws.Range("Y1").Copy() | Out-Null
$ws.Range("Z1:AA1").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Value = "04-19_A"
$ws.Range("AA1").Value = "04-19_0"
$excel.CutCopyMode = 0

# --- Blank rows: Z gets white-fill style (matches X/s=3 pattern); AA stays blank (copy blank format/shape from Y) ---
$ws.Cells.Item(2, 26).Interior.Color = 16777215
$ws.Range("Y2").Copy() | Out-Null
$ws.Range("AA2").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 26).Interior.Color = 16777215
$ws.Range("Y8").Copy() | Out-Null
$ws.Range("AA8").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16, 26).Interior.Color = 16777215
$ws.Range("Y16").Copy() | Out-Null
$ws.Range("AA16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 26).Interior.Color = 16777215
$ws.Range("Y25").Copy() | Out-Null
$ws.Range("AA25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(33, 26).Interior.Color = 16777215
$ws.Range("Y33").Copy() | Out-Null
$ws.Range("AA33").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(50, 26).Interior.Color = 16777215
$ws.Range("Y50").Copy() | Out-Null
$ws.Range("AA50").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(63, 26).Interior.Color = 16777215
$ws.Range("Y63").Copy() | Out-Null
$ws.Range("AA63").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(65, 26).Interior.Color = 16777215
$ws.Range("Y65").Copy() | Out-Null
$ws.Range("AA65").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(75, 26).Interior.Color = 16777215
$ws.Range("Y75").Copy() | Out-Null
$ws.Range("AA75").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(91, 26).Interior.Color = 16777215
$ws.Range("Y91").Copy() | Out-Null
$ws.Range("AA91").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(103, 26).Interior.Color = 16777215
$ws.Range("Y103").Copy() | Out-Null
$ws.Range("AA103").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(104, 26).Interior.Color = 16777215
$ws.Range("Y104").Copy() | Out-Null
$ws.Range("AA104").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(105, 26).Interior.Color = 16777215
$ws.Range("Y105").Copy() | Out-Null
$ws.Range("AA105").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(106, 26).Interior.Color = 16777215
$ws.Range("Y106").Copy() | Out-Null
$ws.Range("AA106").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(107, 26).Interior.Color = 16777215
$ws.Range("Y107").Copy() | Out-Null
$ws.Range("AA107").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(108, 26).Interior.Color = 16777215
$ws.Range("Y108").Copy() | Out-Null
$ws.Range("AA108").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(109, 26).Interior.Color = 16777215
$ws.Range("Y109").Copy() | Out-Null
$ws.Range("AA109").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(110, 26).Interior.Color = 16777215
$ws.Range("Y110").Copy() | Out-Null
$ws.Range("AA110").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111, 26).Interior.Color = 16777215
$ws.Range("Y111").Copy() | Out-Null
$ws.Range("AA111").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(117, 26).Interior.Color = 16777215
$ws.Range("Y117").Copy() | Out-Null
$ws.Range("AA117").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(119, 26).Interior.Color = 16777215
$ws.Range("Y119").Copy() | Out-Null
$ws.Range("AA119").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 66 special case: Y66 numeric; Z66 white-fill blank; AA66 blank (same shape as blank rows) ---
$ws.Cells.Item(66, 26).Interior.Color = 16777215
$ws.Range("Y2").Copy() | Out-Null
$ws.Range("AA66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("Y66").Value = 2524

# --- Data rows: Y -> numeric (same value), Z -> styled numeric, AA -> text formula (converted to static text below) ---
$ws.Range("Y3").Value = 2551
$ws.Cells.Item(3, 26).Interior.Color = 255
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Formula = '="2551"'
$ws.Range("Y4").Value = 0
$ws.Cells.Item(4, 26).Interior.Color = 255
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Formula = '="0"'
$ws.Range("Y5").Value = 4166
$ws.Cells.Item(5, 26).Interior.Color = 255
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Formula = '="4285"'
$ws.Range("Y6").Value = 2529
$ws.Cells.Item(6, 26).Interior.Color = 255
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Formula = '="2522"'
$ws.Range("Y7").Value = 5253
$ws.Cells.Item(7, 26).Interior.Color = 16777215
$ws.Range("Z7").Value = 20
$ws.Range("AA7").Formula = '="5459"'
$ws.Range("Y9").Value = 3635
$ws.Cells.Item(9, 26).Interior.Color = 255
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Formula = '="3650"'
$ws.Range("Y10").Value = 0
$ws.Cells.Item(10, 26).Interior.Color = 255
$ws.Range("Z10").Value = 0
$ws.Range("AA10").Formula = '="0"'
$ws.Range("Y11").Value = 0
$ws.Cells.Item(11, 26).Interior.Color = 255
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Formula = '="0"'
$ws.Range("Y12").Value = 0
$ws.Cells.Item(12, 26).Interior.Color = 255
$ws.Range("Z12").Value = 0
$ws.Range("AA12").Formula = '="0"'
$ws.Range("Y13").Value = 4094
$ws.Cells.Item(13, 26).Interior.Color = 65535
$ws.Range("Z13").Value = 7
$ws.Range("AA13").Formula = '="4028"'
$ws.Range("Y14").Value = 2716
$ws.Cells.Item(14, 26).Interior.Color = 255
$ws.Range("Z14").Value = 0
$ws.Range("AA14").Formula = '="2748"'
$ws.Range("Y15").Value = 4352
$ws.Cells.Item(15, 26).Interior.Color = 16777215
$ws.Range("Z15").Value = 26
$ws.Range("AA15").Formula = '="4444"'
$ws.Range("Y17").Value = 4916
$ws.Cells.Item(17, 26).Interior.Color = 32768
$ws.Range("Z17").Value = 37
$ws.Range("AA17").Formula = '="5003"'
$ws.Range("Y18").Value = 4954
$ws.Cells.Item(18, 26).Interior.Color = 16777215
$ws.Range("Z18").Value = 20
$ws.Range("AA18").Formula = '="5118"'
$ws.Range("Y19").Value = 5599
$ws.Cells.Item(19, 26).Interior.Color = 16777215
$ws.Range("Z19").Value = 23
$ws.Range("AA19").Formula = '="5875"'
$ws.Range("Y20").Value = 4871
$ws.Cells.Item(20, 26).Interior.Color = 32768
$ws.Range("Z20").Value = 34
$ws.Range("AA20").Formula = '="5163"'
$ws.Range("Y21").Value = 5350
$ws.Cells.Item(21, 26).Interior.Color = 16777215
$ws.Range("Z21").Value = 30
$ws.Range("AA21").Formula = '="5604"'
$ws.Range("Y22").Value = 3401
$ws.Cells.Item(22, 26).Interior.Color = 16777215
$ws.Range("Z22").Value = 22
$ws.Range("AA22").Formula = '="3990"'
$ws.Range("Y23").Value = 4259
$ws.Cells.Item(23, 26).Interior.Color = 16777215
$ws.Range("Z23").Value = 20
$ws.Range("AA23").Formula = '="4589"'
$ws.Range("Y24").Value = 0
$ws.Cells.Item(24, 26).Interior.Color = 255
$ws.Range("Z24").Value = 0
$ws.Range("AA24").Formula = '="0"'
$ws.Range("Y26").Value = 4994
$ws.Cells.Item(26, 26).Interior.Color = 32768
$ws.Range("Z26").Value = 33
$ws.Range("AA26").Formula = '="5184"'
$ws.Range("Y27").Value = 3624
$ws.Cells.Item(27, 26).Interior.Color = 16777215
$ws.Range("Z27").Value = 20
$ws.Range("AA27").Formula = '="3944"'
$ws.Range("Y28").Value = 2569
$ws.Cells.Item(28, 26).Interior.Color = 255
$ws.Range("Z28").Value = 0
$ws.Range("AA28").Formula = '="2601"'
$ws.Range("Y29").Value = 4786
$ws.Cells.Item(29, 26).Interior.Color = 16777215
$ws.Range("Z29").Value = 20
$ws.Range("AA29").Formula = '="4945"'
$ws.Range("Y30").Value = 3033
$ws.Cells.Item(30, 26).Interior.Color = 16777215
$ws.Range("Z30").Value = 20
$ws.Range("AA30").Formula = '="3561"'
$ws.Range("Y31").Value = 4715
$ws.Cells.Item(31, 26).Interior.Color = 16777215
$ws.Range("Z31").Value = 30
$ws.Range("AA31").Formula = '="4792"'
$ws.Range("Y32").Value = 4646
$ws.Cells.Item(32, 26).Interior.Color = 16777215
$ws.Range("Z32").Value = 30
$ws.Range("AA32").Formula = '="4636"'
$ws.Range("Y34").Value = 2723
$ws.Cells.Item(34, 26).Interior.Color = 255
$ws.Range("Z34").Value = 0
$ws.Range("AA34").Formula = '="2719"'
$ws.Range("Y35").Value = 4370
$ws.Cells.Item(35, 26).Interior.Color = 65535
$ws.Range("Z35").Value = 3
$ws.Range("AA35").Formula = '="4448"'
$ws.Range("Y36").Value = 3287
$ws.Cells.Item(36, 26).Interior.Color = 255
$ws.Range("Z36").Value = 0
$ws.Range("AA36").Formula = '="3351"'
$ws.Range("Y37").Value = 5306
$ws.Cells.Item(37, 26).Interior.Color = 16777215
$ws.Range("Z37").Value = 30
$ws.Range("AA37").Formula = '="5472"'
$ws.Range("Y38").Value = 4123
$ws.Cells.Item(38, 26).Interior.Color = 16777215
$ws.Range("Z38").Value = 30
$ws.Range("AA38").Formula = '="4062"'
$ws.Range("Y39").Value = 3515
$ws.Cells.Item(39, 26).Interior.Color = 65535
$ws.Range("Z39").Value = 11
$ws.Range("AA39").Formula = '="3703"'
$ws.Range("Y40").Value = 5395
$ws.Cells.Item(40, 26).Interior.Color = 32768
$ws.Range("Z40").Value = 33
$ws.Range("AA40").Formula = '="5488"'
$ws.Range("Y41").Value = 5201
$ws.Cells.Item(41, 26).Interior.Color = 16777215
$ws.Range("Z41").Value = 30
$ws.Range("AA41").Formula = '="5329"'
$ws.Range("Y42").Value = 4789
$ws.Cells.Item(42, 26).Interior.Color = 16777215
$ws.Range("Z42").Value = 30
$ws.Range("AA42").Formula = '="4892"'
$ws.Range("Y43").Value = 4881
$ws.Cells.Item(43, 26).Interior.Color = 16777215
$ws.Range("Z43").Value = 23
$ws.Range("AA43").Formula = '="5002"'
$ws.Range("Y44").Value = 4654
$ws.Cells.Item(44, 26).Interior.Color = 16777215
$ws.Range("Z44").Value = 20
$ws.Range("AA44").Formula = '="4868"'
$ws.Range("Y45").Value = 4944
$ws.Cells.Item(45, 26).Interior.Color = 16777215
$ws.Range("Z45").Value = 30
$ws.Range("AA45").Formula = '="5010"'
$ws.Range("Y46").Value = 4738
$ws.Cells.Item(46, 26).Interior.Color = 16777215
$ws.Range("Z46").Value = 20
$ws.Range("AA46").Formula = '="4861"'
$ws.Range("Y47").Value = 4318
$ws.Cells.Item(47, 26).Interior.Color = 65535
$ws.Range("Z47").Value = 16
$ws.Range("AA47").Formula = '="4580"'
$ws.Range("Y48").Value = 5034
$ws.Cells.Item(48, 26).Interior.Color = 16777215
$ws.Range("Z48").Value = 20
$ws.Range("AA48").Formula = '="5262"'
$ws.Range("Y49").Value = 4195
$ws.Cells.Item(49, 26).Interior.Color = 65535
$ws.Range("Z49").Value = 15
$ws.Range("AA49").Formula = '="4327"'
$ws.Range("Y51").Value = 4170
$ws.Cells.Item(51, 26).Interior.Color = 16777215
$ws.Range("Z51").Value = 23
$ws.Range("AA51").Formula = '="4276"'
$ws.Range("Y52").Value = 4231
$ws.Cells.Item(52, 26).Interior.Color = 16777215
$ws.Range("Z52").Value = 20
$ws.Range("AA52").Formula = '="4310"'
$ws.Range("Y53").Value = 4213
$ws.Cells.Item(53, 26).Interior.Color = 16777215
$ws.Range("Z53").Value = 20
$ws.Range("AA53").Formula = '="4302"'
$ws.Range("Y54").Value = 4029
$ws.Cells.Item(54, 26).Interior.Color = 16777215
$ws.Range("Z54").Value = 30
$ws.Range("AA54").Formula = '="4108"'
$ws.Range("Y55").Value = 3998
$ws.Cells.Item(55, 26).Interior.Color = 16777215
$ws.Range("Z55").Value = 30
$ws.Range("AA55").Formula = '="4100"'
$ws.Range("Y56").Value = 4035
$ws.Cells.Item(56, 26).Interior.Color = 16777215
$ws.Range("Z56").Value = 22
$ws.Range("AA56").Formula = '="4094"'
$ws.Range("Y57").Value = 4261
$ws.Cells.Item(57, 26).Interior.Color = 16777215
$ws.Range("Z57").Value = 20
$ws.Range("AA57").Formula = '="4368"'
$ws.Range("Y58").Value = 0
$ws.Cells.Item(58, 26).Interior.Color = 255
$ws.Range("Z58").Value = 0
$ws.Range("AA58").Formula = '="0"'
$ws.Range("Y59").Value = 2777
$ws.Cells.Item(59, 26).Interior.Color = 255
$ws.Range("Z59").Value = 0
$ws.Range("AA59").Formula = '="2851"'
$ws.Range("Y60").Value = 2465
$ws.Cells.Item(60, 26).Interior.Color = 255
$ws.Range("Z60").Value = 0
$ws.Range("AA60").Formula = '="2461"'
$ws.Range("Y61").Value = 4026
$ws.Cells.Item(61, 26).Interior.Color = 255
$ws.Range("Z61").Value = 0
$ws.Range("AA61").Formula = '="4024"'
$ws.Range("Y62").Value = 1987
$ws.Cells.Item(62, 26).Interior.Color = 255
$ws.Range("Z62").Value = 0
$ws.Range("AA62").Formula = '="1981"'
$ws.Range("Y64").Value = 0
$ws.Cells.Item(64, 26).Interior.Color = 255
$ws.Range("Z64").Value = 0
$ws.Range("AA64").Formula = '="0"'
$ws.Range("Y67").Value = 5173
$ws.Cells.Item(67, 26).Interior.Color = 32768
$ws.Range("Z67").Value = 33
$ws.Range("AA67").Formula = '="5295"'
$ws.Range("Y68").Value = 0
$ws.Cells.Item(68, 26).Interior.Color = 255
$ws.Range("Z68").Value = 0
$ws.Range("AA68").Formula = '="0"'
$ws.Range("Y69").Value = 2618
$ws.Cells.Item(69, 26).Interior.Color = 255
$ws.Range("Z69").Value = 0
$ws.Range("AA69").Formula = '="2660"'
$ws.Range("Y70").Value = 0
$ws.Cells.Item(70, 26).Interior.Color = 255
$ws.Range("Z70").Value = 0
$ws.Range("AA70").Formula = '="0"'
$ws.Range("Y71").Value = 3992
$ws.Cells.Item(71, 26).Interior.Color = 16777215
$ws.Range("Z71").Value = 20
$ws.Range("AA71").Formula = '="4347"'
$ws.Range("Y72").Value = 3234
$ws.Cells.Item(72, 26).Interior.Color = 255
$ws.Range("Z72").Value = 0
$ws.Range("AA72").Formula = '="3271"'
$ws.Range("Y73").Value = 0
$ws.Cells.Item(73, 26).Interior.Color = 255
$ws.Range("Z73").Value = 0
$ws.Range("AA73").Formula = '="0"'
$ws.Range("Y74").Value = 0
$ws.Cells.Item(74, 26).Interior.Color = 255
$ws.Range("Z74").Value = 0
$ws.Range("AA74").Formula = '="1229"'
$ws.Range("Y76").Value = 0
$ws.Cells.Item(76, 26).Interior.Color = 255
$ws.Range("Z76").Value = 0
$ws.Range("AA76").Formula = '="0"'
$ws.Range("Y77").Value = 3019
$ws.Cells.Item(77, 26).Interior.Color = 255
$ws.Range("Z77").Value = 0
$ws.Range("AA77").Formula = '="3043"'
$ws.Range("Y78").Value = 0
$ws.Cells.Item(78, 26).Interior.Color = 255
$ws.Range("Z78").Value = 0
$ws.Range("AA78").Formula = '="0"'
$ws.Range("Y79").Value = 0
$ws.Cells.Item(79, 26).Interior.Color = 255
$ws.Range("Z79").Value = 0
$ws.Range("AA79").Formula = '="0"'
$ws.Range("Y80").Value = 1490
$ws.Cells.Item(80, 26).Interior.Color = 255
$ws.Range("Z80").Value = 0
$ws.Range("AA80").Formula = '="1489"'
$ws.Range("Y81").Value = 0
$ws.Cells.Item(81, 26).Interior.Color = 255
$ws.Range("Z81").Value = 0
$ws.Range("AA81").Formula = '="0"'
$ws.Range("Y82").Value = 0
$ws.Cells.Item(82, 26).Interior.Color = 255
$ws.Range("Z82").Value = 0
$ws.Range("AA82").Formula = '="0"'
$ws.Range("Y83").Value = 0
$ws.Cells.Item(83, 26).Interior.Color = 255
$ws.Range("Z83").Value = 0
$ws.Range("AA83").Formula = '="0"'
$ws.Range("Y84").Value = 0
$ws.Cells.Item(84, 26).Interior.Color = 255
$ws.Range("Z84").Value = 0
$ws.Range("AA84").Formula = '="0"'
$ws.Range("Y85").Value = 0
$ws.Cells.Item(85, 26).Interior.Color = 255
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Formula = '="0"'
$ws.Range("Y86").Value = 0
$ws.Cells.Item(86, 26).Interior.Color = 255
$ws.Range("Z86").Value = 0
$ws.Range("AA86").Formula = '="0"'
$ws.Range("Y87").Value = 0
$ws.Cells.Item(87, 26).Interior.Color = 255
$ws.Range("Z87").Value = 0
$ws.Range("AA87").Formula = '="0"'
$ws.Range("Y88").Value = 2659
$ws.Cells.Item(88, 26).Interior.Color = 255
$ws.Range("Z88").Value = 0
$ws.Range("AA88").Formula = '="2746"'
$ws.Range("Y89").Value = 0
$ws.Cells.Item(89, 26).Interior.Color = 255
$ws.Range("Z89").Value = 0
$ws.Range("AA89").Formula = '="0"'
$ws.Range("Y90").Value = 0
$ws.Cells.Item(90, 26).Interior.Color = 255
$ws.Range("Z90").Value = 0
$ws.Range("AA90").Formula = '="0"'
$ws.Range("Y92").Value = 0
$ws.Cells.Item(92, 26).Interior.Color = 255
$ws.Range("Z92").Value = 0
$ws.Range("AA92").Formula = '="0"'
$ws.Range("Y93").Value = 3438
$ws.Cells.Item(93, 26).Interior.Color = 16777215
$ws.Range("Z93").Value = 20
$ws.Range("AA93").Formula = '="3753"'
$ws.Range("Y94").Value = 0
$ws.Cells.Item(94, 26).Interior.Color = 255
$ws.Range("Z94").Value = 0
$ws.Range("AA94").Formula = '="0"'
$ws.Range("Y95").Value = 0
$ws.Cells.Item(95, 26).Interior.Color = 255
$ws.Range("Z95").Value = 0
$ws.Range("AA95").Formula = '="0"'
$ws.Range("Y96").Value = 0
$ws.Cells.Item(96, 26).Interior.Color = 255
$ws.Range("Z96").Value = 0
$ws.Range("AA96").Formula = '="0"'
$ws.Range("Y97").Value = 0
$ws.Cells.Item(97, 26).Interior.Color = 255
$ws.Range("Z97").Value = 0
$ws.Range("AA97").Formula = '="0"'
$ws.Range("Y98").Value = 0
$ws.Cells.Item(98, 26).Interior.Color = 255
$ws.Range("Z98").Value = 0
$ws.Range("AA98").Formula = '="0"'
$ws.Range("Y99").Value = 0
$ws.Cells.Item(99, 26).Interior.Color = 255
$ws.Range("Z99").Value = 0
$ws.Range("AA99").Formula = '="0"'
$ws.Range("Y100").Value = 0
$ws.Cells.Item(100, 26).Interior.Color = 255
$ws.Range("Z100").Value = 0
$ws.Range("AA100").Formula = '="0"'
$ws.Range("Y101").Value = 0
$ws.Cells.Item(101, 26).Interior.Color = 255
$ws.Range("Z101").Value = 0
$ws.Range("AA101").Formula = '="0"'
$ws.Range("Y102").Value = 0
$ws.Cells.Item(102, 26).Interior.Color = 255
$ws.Range("Z102").Value = 0
$ws.Range("AA102").Formula = '="0"'
$ws.Range("Y112").Value = 4364
$ws.Cells.Item(112, 26).Interior.Color = 255
$ws.Range("Z112").Value = 0
$ws.Range("AA112").Formula = '="4364"'
$ws.Range("Y113").Value = 2819
$ws.Cells.Item(113, 26).Interior.Color = 255
$ws.Range("Z113").Value = 0
$ws.Range("AA113").Formula = '="2898"'
$ws.Range("Y114").Value = 1587
$ws.Cells.Item(114, 26).Interior.Color = 255
$ws.Range("Z114").Value = 0
$ws.Range("AA114").Formula = '="1587"'
$ws.Range("Y115").Value = 5050
$ws.Cells.Item(115, 26).Interior.Color = 16777215
$ws.Range("Z115").Value = 23
$ws.Range("AA115").Formula = '="5209"'
$ws.Range("Y116").Value = 3654
$ws.Cells.Item(116, 26).Interior.Color = 16777215
$ws.Range("Z116").Value = 23
$ws.Range("AA116").Formula = '="3704"'
$ws.Range("Y118").Value = 2829
$ws.Cells.Item(118, 26).Interior.Color = 255
$ws.Range("Z118").Value = 0
$ws.Range("AA118").Formula = '="2815"'
$ws.Range("Y120").Value = 1464
$ws.Cells.Item(120, 26).Interior.Color = 65535
$ws.Range("Z120").Value = 1
$ws.Range("AA120").Formula = '="1477"'
$ws.Range("Y121").Value = 0
$ws.Cells.Item(121, 26).Interior.Color = 255
$ws.Range("Z121").Value = 0
$ws.Range("AA121").Formula = '="0"'
$ws.Range("Y122").Value = 2626
$ws.Cells.Item(122, 26).Interior.Color = 255
$ws.Range("Z122").Value = 0
$ws.Range("AA122").Formula = '="2723"'
$ws.Range("Y123").Value = 2578
$ws.Cells.Item(123, 26).Interior.Color = 255
$ws.Range("Z123").Value = 0
$ws.Range("AA123").Formula = '="2576"'
$ws.Range("Y124").Value = 1441
$ws.Cells.Item(124, 26).Interior.Color = 255
$ws.Range("Z124").Value = 0
$ws.Range("AA124").Formula = '="1445"'
$ws.Range("Y125").Value = 2534
$ws.Cells.Item(125, 26).Interior.Color = 255
$ws.Range("Z125").Value = 0
$ws.Range("AA125").Formula = '="2584"'
$ws.Range("Y126").Value = 1506
$ws.Cells.Item(126, 26).Interior.Color = 65535
$ws.Range("Z126").Value = 10
$ws.Range("AA126").Formula = '="1674"'
$ws.Range("Y127").Value = 1490
$ws.Cells.Item(127, 26).Interior.Color = 65535
$ws.Range("Z127").Value = 2
$ws.Range("AA127").Formula = '="1537"'

# Convert AA-column formulas to static text values (keeps string type, no style churn)
$ws.Range("AA3:AA127").Copy() | Out-Null
$ws.Range("AA3:AA127").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0